$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.758.55'
$ws.Range("E2").Value = '  -2.38%  '

$ws.Range("D3").Value = '2.499.16'
$ws.Range("E3").Value = '  -4.82%  '

$ws.Range("E4").Value = '  +0.00%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Formula = '576.98'
$c.ClearFormats()
$ws.Range("E5").Value = '  -3.02%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Formula = '166.69'
$c.ClearFormats()
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D9").Value = '2.498.47'
$ws.Range("E9").Value = '  -4.83%  '

$ws.Range("E10").Value = '  -1.82%  '

$ws.Range("E11").Value = '  -0.28%  '

$ws.Range("E12").Value = '  -4.52%  '

$ws.Range("E13").Value = '  -3.00%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Formula = '26.15'
$c.ClearFormats()
$ws.Range("E14").Value = '  -5.47%  '

$ws.Range("D15").Value = '2.957.97'
$ws.Range("E15").Value = '  -4.68%  '

$ws.Range("E16").Value = '  -4.68%  '

$ws.Range("D17").Value = '65.578.73'
$ws.Range("E17").Value = '  -2.29%  '

$ws.Range("D18").Value = '2.502.64'
$ws.Range("E18").Value = '  -4.50%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Formula = '11.19'
$c.ClearFormats()
$ws.Range("E19").Value = '  -7.30%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Formula = '7.57'
$c.ClearFormats()
$ws.Range("E20").Value = '  -5.01%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Formula = '342.95'
$c.ClearFormats()
$ws.Range("E21").Value = '  -4.12%  '

$ws.Range("E22").Value = '  -3.40%  '

$ws.Range("E23").Value = '  -3.16%  '

$ws.Range("E24").Value = '  -0.04%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Formula = '1.93'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.51%  '

$ws.Range("E26").Value = '  -1.37%  '

$ws.Range("E27").Value = '  -3.72%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Formula = '0.999'
$c.ClearFormats()
$ws.Range("E28").Value = '  -0.20%  '

$ws.Range("D29").Value = '2.629.97'
$ws.Range("E29").Value = '  -4.70%  '

$ws.Range("E30").Value = '  -3.34%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Formula = '8.09'
$c.ClearFormats()
$ws.Range("E31").Value = '  +2.03%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Formula = '518.68'
$c.ClearFormats()
$ws.Range("E32").Value = '  -5.04%  '

$ws.Range("E33").Value = '  -3.48%  '

$ws.Range("E34").Value = '  -5.41%  '

$ws.Range("E35").Value = '  -4.35%  '

$ws.Range("E36").Value = '  -0.04%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Formula = '157.07'
$c.ClearFormats()
$ws.Range("E37").Value = '  +0.24%  '

$ws.Range("E38").Value = '  -4.54%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Formula = '18.43'
$c.ClearFormats()
$ws.Range("E39").Value = '  -3.06%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Formula = '18.28'
$c.ClearFormats()
$ws.Range("E40").Value = '  +0.51%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Formula = '0.352'
$c.ClearFormats()
$ws.Range("E41").Value = '  -3.84%  '

$ws.Range("E42").Value = '  -3.56%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Formula = '5.01'
$c.ClearFormats()
$ws.Range("E43").Value = '  -3.92%  '

$ws.Range("E44").Value = '  -0.06%  '

$ws.Range("E45").Value = '  +0.04%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Formula = '145.95'
$c.ClearFormats()
$ws.Range("E46").Value = '  -4.20%  '

$ws.Range("E47").Value = '  -4.95%  '

$ws.Range("E48").Value = '  -2.88%  '

$ws.Range("D49").Value = '0.0₆0272'
$ws.Range("E49").Value = '  -8.33%  '

$ws.Range("E50").Value = '  +0.03%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Formula = '0.0750'
$c.ClearFormats()
$ws.Range("E51").Value = '  -2.77%  '

